# negyedik dia címének megadása
# (1) Merge the split runs in slide 4's third content paragraph back into
#     a single run - matches the authored edit where the three runs
#     "Legyünk " / "képesek külön-külön " / "szerkeszteni ..." become one.
$p = $ppt.ActivePresentation

$s4 = $p.Slides.Item(4)
$content4 = $s4.Shapes.Item(2).TextFrame.TextRange
$para3 = $content4.Paragraphs(3)
$para3.Text = "placeholder"
$para3b = $content4.Paragraphs(3)
$para3b.Text = "Legyünk képesek külön-külön szerkeszteni az ilyen jellegű autókat (melyek egy tulajdonoshoz vannak rendelve.)"

# (2) Add a new 5th slide ("A szoftver jelenlegi helyzete") using the same
#     "Title and Content" layout as the rest of the deck. Duplicating the
#     existing slide 4 (rather than Slides.Add) keeps the slide's internal
#     structure (xfrm, extLst/creationId, clrMapOvr, placeholder names)
#     consistent with the rest of the presentation.
$dup = $s4.Duplicate()
$s5 = $dup.Item(1)

$title5 = $s5.Shapes.Item(1).TextFrame.TextRange
$title5.Text = "A szoftver jelenlegi helyzete"

# Clear out the copied body text paragraphs, leaving a single empty
# paragraph behind (the new slide has no body content yet).
$content5 = $s5.Shapes.Item(2).TextFrame.TextRange
for ($i = 0; $i -lt 4; $i++) {
    $content5.Paragraphs(1).Delete()
}
